$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1975.32
$ws.Range("F2").Value = 147.49
$ws.Range("G2").Value = 1273.05
$ws.Range("K2").Value = 1
$ws.Range("M2").Value = 1
$ws.Range("N2").Value = 0

# Row 3
$ws.Range("E3").Value = 72.18000000000001
$ws.Range("F3").Value = 1737.21
$ws.Range("G3").Value = 523.89
$ws.Range("K3").Value = 1
$ws.Range("M3").Value = 0.34
$ws.Range("N3").Value = 0.99

# Row 4
$ws.Range("E4").Value = 67.22
$ws.Range("F4").Value = 631.79
$ws.Range("G4").Value = 430.11
$ws.Range("K4").Value = 1
$ws.Range("M4").Value = 1.27
$ws.Range("N4").Value = 0

# Row 5
$ws.Range("E5").Value = 39.18
$ws.Range("F5").Value = 424.69
$ws.Range("G5").Value = 92738.05
$ws.Range("K5").Value = 1
$ws.Range("M5").Value = 1702.06
$ws.Range("N5").Value = 11.05
$ws.Range("O5").Value = 49.7

# Row 6
$ws.Range("E6").Value = 4.57
$ws.Range("F6").Value = 3.19
$ws.Range("G6").Value = 0
$ws.Range("K6").Value = 1
$ws.Range("M6").Value = 176791.47
$ws.Range("N6").Value = 332
$ws.Range("O6").Value = 14191.81
